# "task정의 및 일정.xlsx" — update 개발목록 progress (%) figures and
# restore the author's last cursor/scroll position before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("개발목록")
$ws.Activate()

# Progress column (I) updates
$ws.Range("I7").Value  = 0.6   # 0.5 -> 0.6
$ws.Range("I9").Value  = 0.9   # 0.8 -> 0.9
$ws.Range("I10").Value = 0.9   # 0.5 -> 0.9
$ws.Range("I11").Value = 0.5   # 0.4 -> 0.5
$ws.Range("I13").Value = 0.5   # 0.4 -> 0.5
$ws.Range("I14").Value = 0.5   # 0.4 -> 0.5
$ws.Range("I15").Value = 0.5   # 0.4 -> 0.5
$ws.Range("I16").Value = 0.5   # 0.4 -> 0.5

# Restore the saved selection/scroll state (frozen pane top-left row moves
# from A6 to A3, active cell moves from I17 to I11)
$ws.Range("I11").Select()
$ws.Application.ActiveWindow.ScrollRow = 3
